$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Sema6a","Plxna4","ECs",3,1,53.11151633333333,159.334549,0.7183308303902657,0.7183308303902658,3,1,3.316186,9.948558,0.7638813129544791,0.7638813129544793,176.1276669033713,1585.149002130342,0.5487194978541974,0.5487194978541976),
    @("ECs","Sema6a","Plxna4","FAPs",3,1,53.11151633333333,159.334549,0.7183308303902657,0.7183308303902658,3,1,0.4282866666666667,1.28486,0.09865555829927233,0.09865555829927235,22.74695429201556,204.72258862814,0.07086732911573156,0.07086732911573158),
    @("ECs","Sema6a","Plxna4","M2",3,1,53.11151633333333,159.334549,0.7183308303902657,0.7183308303902658,3,1,0.417289,1.251867,0.09612225285356782,0.09612225285356783,22.16285153922033,199.465663852983,0.06904757771128646,0.06904757771128647),
    @("ECs","Sema6a","Plxna4","sCs",3,1,53.11151633333333,159.334549,0.7183308303902657,0.7183308303902658,3,1,0.1794703333333333,0.538411,0.04134087589268053,0.04134087589268053,9.53194154018211,85.787473861639,0.02969642570905012,0.02969642570905013),
    @("FAPs","Sema6a","Plxna4","ECs",3,1,4.187132,12.561396,0.05663076888328195,0.05663076888328195,3,1,3.316186,9.948558,0.7638813129544791,0.7638813129544793,13.885308518552,124.967776666968,0.04325918608818308,0.04325918608818308),
    @("FAPs","Sema6a","Plxna4","FAPs",3,1,4.187132,12.561396,0.05663076888328195,0.05663076888328195,3,1,0.4282866666666667,1.28486,0.09865555829927233,0.09865555829927235,1.793292807173334,16.13963526456,0.005586940121097239,0.00558694012109724),
    @("FAPs","Sema6a","Plxna4","M2",3,1,4.187132,12.561396,0.05663076888328195,0.05663076888328195,3,1,0.417289,1.251867,0.09612225285356782,0.09612225285356783,1.747244125148,15.725197126332,0.005443477085890788,0.005443477085890789),
    @("FAPs","Sema6a","Plxna4","sCs",3,1,4.187132,12.561396,0.05663076888328195,0.05663076888328195,3,1,0.1794703333333333,0.538411,0.04134087589268053,0.04134087589268053,0.7514659757506666,6.763193781756,0.002341165588110834,0.002341165588110834),
    @("M2","Sema6a","Plxna4","ECs",1,0.3333333333333333,0.02264733333333334,0.067942,0.0003063041479997877,0.0003063041479997877,3,1,3.316186,9.948558,0.7638813129544791,0.7638813129544793,0.07510276973733335,0.6759249276360001,0.000233980014737481,0.000233980014737481),
    @("M2","Sema6a","Plxna4","FAPs",1,0.3333333333333333,0.02264733333333334,0.067942,0.0003063041479997877,0.0003063041479997877,3,1,0.4282866666666667,1.28486,0.09865555829927233,0.09865555829927235,0.009699550902222224,0.08729595812000002,0.000030218606730302,0.000030218606730302),
    @("M2","Sema6a","Plxna4","M2",1,0.3333333333333333,0.02264733333333334,0.067942,0.0003063041479997877,0.0003063041479997877,3,1,0.417289,1.251867,0.09612225285356782,0.09612225285356783,0.009450483079333334,0.085054347714,0.02097996460008708,0.02097996460008708),
    @("M2","Sema6a","Plxna4","sCs",1,0.3333333333333333,0.02264733333333334,0.067942,0.0003063041479997877,0.0003063041479997877,3,1,0.1794703333333333,0.538411,0.04134087589268053,0.04134087589268053,0.004064524462444445,0.036580720162,0.00001266288176787248,0.00001266288176787248),
    @("sCs","Sema6a","Plxna4","ECs",3,1,16.61610766666667,49.848323,0.2247320965784526,0.2247320965784526,3,1,3.316186,9.948558,0.7638813129544791,0.7638813129544793,55.10210361869267,495.918932568234,0.1716686489973611,0.1716686489973612),
    @("sCs","Sema6a","Plxna4","FAPs",3,1,16.61610766666667,49.848323,0.2247320965784526,0.2247320965784526,3,1,0.4282866666666667,1.28486,0.09865555829927233,0.09865555829927235,7.116457365531113,64.04811628978001,0.02217107045571323,0.02217107045571323),
    @("sCs","Sema6a","Plxna4","M2",3,1,16.61610766666667,49.848323,0.2247320965784526,0.2247320965784526,3,1,0.417289,1.251867,0.09612225285356782,0.09612225285356783,6.933718952115668,62.403470569041,0.02160175541162644,0.02160175541162644),
    @("sCs","Sema6a","Plxna4","sCs",3,1,16.61610766666667,49.848323,0.2247320965784526,0.2247320965784526,3,1,0.1794703333333333,0.538411,0.04134087589268053,0.04134087589268053,2.982098381639222,26.838885434753,0.009290621713751702,0.009290621713751702)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Output "Updated $($data.Count) rows"
